# Remove the RG (ID card) reference from the certificate text, keeping
# only the CPF reference, and tidy up a couple of neighbouring runs that
# get merged together as a side effect of re-flowing the paragraph.
#
# Before: "... {{NOME}}, portador do RG nº {{RG}} e CPF nº {{CPF}}, concluiu ..."
# After:  "... {{NOME}}, portador do CPF nº {{CPF}}, concluiu ..."
#
# Also: "São Carlos" + ", " + "{{DATA}}"  ->  "São Carlos, " + "{{DATA}}"
#
# Target shape: slide 1, shape "Rectangle 5" (the certificate body text).

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tf  = $shp.TextFrame
$tr  = $tf.TextRange

# ---------------------------------------------------------------------
# Paragraph 1 - locate the relevant words/placeholders by searching the
# paragraph text so the edit does not depend on hard-coded offsets.
# Character positions are then patched right-to-left (highest offset
# first) so earlier offsets stay valid while later text is resized.
# ---------------------------------------------------------------------
$para1 = $tr.Paragraphs(1, 1)
$full1 = $para1.Text

$idxPortador = $full1.IndexOf("portador do RG")      # start of "portador do RG nº "
$idxRG       = $full1.IndexOf("{{RG}}")               # start of "{{RG}} "
$idxCPFTag   = $full1.IndexOf("{{CPF}}")               # start of "{{CPF}}"
$idxConcluiu = $full1.IndexOf("concluiu")              # start of "concluiu..." word

# 1-based character positions (PowerPoint Characters()/Start are 1-based)
$posPortador = $idxPortador + 1
$posRG       = $idxRG + 1
$posCPFTag   = $idxCPFTag + 1
$posConcluiu = $idxConcluiu + 1

# "e CPF nº " runs from just after "{{RG}} " up to just before "{{CPF}}"
$posECPF   = $posRG + "{{RG}} ".Length
$lenECPF   = $posCPFTag - $posECPF

# the ", " between "{{CPF}}" and "concluiu" sits right before "concluiu"
$posComma  = $posCPFTag + "{{CPF}}".Length
$lenComma  = $posConcluiu - $posComma

# 1) Drop the ", " that used to separate "{{CPF}}" from "concluiu..." -
#    it gets folded into the start of the "concluiu..." run instead.
$comma = $para1.Characters($posComma, $lenComma)
$comma.Text = ""
$concluiuRun = $para1.Characters($posComma, ("concluiu com aproveitamento o curso da ").Length)
$concluiuRun.Text = ", concluiu com aproveitamento o curso da "

# 2) "e CPF nº " -> "CPF nº " (drop the leading "e ")
$ePrefix = $para1.Characters($posECPF, 2)
$ePrefix.Text = ""

# 3) Split "CPF nº " into "CPF " + "nº " (replacing just the leading
#    sub-range splits the run in two, keeping the original formatting).
$cpfWord = $para1.Characters($posECPF, "CPF ".Length)
$cpfWord.Text = "CPF "

# 4) "{{RG}} " -> "do "
$rgRun = $para1.Characters($posRG, "{{RG}} ".Length)
$rgRun.Text = "do "

# 5) "portador do RG nº " -> "portador "
$portadorRun = $para1.Characters($posPortador, "portador do RG nº ".Length)
$portadorRun.Text = "portador "

# ---------------------------------------------------------------------
# Paragraph 4: "São Carlos" + ", " + "{{DATA}}"  ->  "São Carlos, " + "{{DATA}}"
# ---------------------------------------------------------------------
$para4 = $tr.Paragraphs(4, 1)
$full4 = $para4.Text

$idxDataTag = $full4.IndexOf("{{DATA}}")
$posCidade  = 1
$lenCidade  = "São Carlos".Length
$posSep     = $posCidade + $lenCidade
$lenSep     = $idxDataTag - ($posSep - 1)

$sep = $para4.Characters($posSep, $lenSep)
$sep.Text = ""

$cidade = $para4.Characters($posCidade, $lenCidade)
$cidade.Text = "São Carlos, "
